$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values per row
# following re-running SGNN dialog act annotation after transcript cleanup.

$ws.Cells.Item(4, 9).Value = "sd"
$ws.Cells.Item(4, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(8, 9).Value = "sd"
$ws.Cells.Item(8, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(23, 9).Value = "%"
$ws.Cells.Item(23, 10).Value = "Uninterpretable"
$ws.Cells.Item(24, 9).Value = "sd"
$ws.Cells.Item(24, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(26, 9).Value = "%"
$ws.Cells.Item(26, 10).Value = "Uninterpretable"
$ws.Cells.Item(28, 9).Value = "aa"
$ws.Cells.Item(28, 10).Value = "Agree/Accept"
$ws.Cells.Item(51, 9).Value = "sd"
$ws.Cells.Item(51, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(68, 9).Value = "sd"
$ws.Cells.Item(68, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(69, 9).Value = "aa"
$ws.Cells.Item(69, 10).Value = "Agree/Accept"
$ws.Cells.Item(79, 9).Value = "ba"
$ws.Cells.Item(79, 10).Value = "Appreciation"
$ws.Cells.Item(81, 9).Value = "sd"
$ws.Cells.Item(81, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(83, 9).Value = "sv"
$ws.Cells.Item(83, 10).Value = "Statement-opinion"
$ws.Cells.Item(86, 9).Value = "sd"
$ws.Cells.Item(86, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(98, 9).Value = "ba"
$ws.Cells.Item(98, 10).Value = "Appreciation"
$ws.Cells.Item(111, 9).Value = "%"
$ws.Cells.Item(111, 10).Value = "Uninterpretable"
$ws.Cells.Item(136, 9).Value = "aa"
$ws.Cells.Item(136, 10).Value = "Agree/Accept"
$ws.Cells.Item(153, 9).Value = "aa"
$ws.Cells.Item(153, 10).Value = "Agree/Accept"
$ws.Cells.Item(154, 9).Value = "%"
$ws.Cells.Item(154, 10).Value = "Uninterpretable"
$ws.Cells.Item(155, 9).Value = "aa"
$ws.Cells.Item(155, 10).Value = "Agree/Accept"
$ws.Cells.Item(158, 9).Value = "aa"
$ws.Cells.Item(158, 10).Value = "Agree/Accept"
$ws.Cells.Item(159, 9).Value = "aa"
$ws.Cells.Item(159, 10).Value = "Agree/Accept"
$ws.Cells.Item(160, 9).Value = "aa"
$ws.Cells.Item(160, 10).Value = "Agree/Accept"
$ws.Cells.Item(161, 9).Value = "aa"
$ws.Cells.Item(161, 10).Value = "Agree/Accept"
$ws.Cells.Item(170, 9).Value = "%"
$ws.Cells.Item(170, 10).Value = "Uninterpretable"
$ws.Cells.Item(171, 9).Value = "%"
$ws.Cells.Item(171, 10).Value = "Uninterpretable"
$ws.Cells.Item(177, 9).Value = "aa"
$ws.Cells.Item(177, 10).Value = "Agree/Accept"
$ws.Cells.Item(178, 9).Value = "aa"
$ws.Cells.Item(178, 10).Value = "Agree/Accept"
$ws.Cells.Item(179, 9).Value = "aa"
$ws.Cells.Item(179, 10).Value = "Agree/Accept"
$ws.Cells.Item(180, 9).Value = "aa"
$ws.Cells.Item(180, 10).Value = "Agree/Accept"
$ws.Cells.Item(191, 9).Value = "%"
$ws.Cells.Item(191, 10).Value = "Uninterpretable"
$ws.Cells.Item(192, 9).Value = "%"
$ws.Cells.Item(192, 10).Value = "Uninterpretable"
$ws.Cells.Item(201, 9).Value = "sd"
$ws.Cells.Item(201, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(204, 9).Value = "aa"
$ws.Cells.Item(204, 10).Value = "Agree/Accept"
$ws.Cells.Item(206, 9).Value = "%"
$ws.Cells.Item(206, 10).Value = "Uninterpretable"
$ws.Cells.Item(223, 9).Value = "sd"
$ws.Cells.Item(223, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(226, 9).Value = "sv"
$ws.Cells.Item(226, 10).Value = "Statement-opinion"
$ws.Cells.Item(227, 9).Value = "aa"
$ws.Cells.Item(227, 10).Value = "Agree/Accept"
$ws.Cells.Item(232, 9).Value = "aa"
$ws.Cells.Item(232, 10).Value = "Agree/Accept"
$ws.Cells.Item(234, 9).Value = "sd"
$ws.Cells.Item(234, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(235, 9).Value = "%"
$ws.Cells.Item(235, 10).Value = "Uninterpretable"
$ws.Cells.Item(249, 9).Value = "sd"
$ws.Cells.Item(249, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(252, 9).Value = "aa"
$ws.Cells.Item(252, 10).Value = "Agree/Accept"
$ws.Cells.Item(264, 9).Value = "aa"
$ws.Cells.Item(264, 10).Value = "Agree/Accept"
$ws.Cells.Item(265, 9).Value = "aa"
$ws.Cells.Item(265, 10).Value = "Agree/Accept"
$ws.Cells.Item(269, 9).Value = "sd"
$ws.Cells.Item(269, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(270, 9).Value = "sd"
$ws.Cells.Item(270, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(271, 9).Value = "%"
$ws.Cells.Item(271, 10).Value = "Uninterpretable"
$ws.Cells.Item(272, 9).Value = "sd"
$ws.Cells.Item(272, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(274, 9).Value = "%"
$ws.Cells.Item(274, 10).Value = "Uninterpretable"
$ws.Cells.Item(277, 9).Value = "aa"
$ws.Cells.Item(277, 10).Value = "Agree/Accept"
$ws.Cells.Item(278, 9).Value = "aa"
$ws.Cells.Item(278, 10).Value = "Agree/Accept"
$ws.Cells.Item(290, 9).Value = "aa"
$ws.Cells.Item(290, 10).Value = "Agree/Accept"
$ws.Cells.Item(292, 9).Value = "sd"
$ws.Cells.Item(292, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(295, 9).Value = "sd"
$ws.Cells.Item(295, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(298, 9).Value = "sd"
$ws.Cells.Item(298, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(307, 9).Value = "%"
$ws.Cells.Item(307, 10).Value = "Uninterpretable"
$ws.Cells.Item(313, 9).Value = "sd"
$ws.Cells.Item(313, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(333, 9).Value = "sd"
$ws.Cells.Item(333, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(339, 9).Value = "sd"
$ws.Cells.Item(339, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(341, 9).Value = "sd"
$ws.Cells.Item(341, 10).Value = "Statement-non-opinion"
